# Denmark Division 1 — base update (2024-01-29 17:07)
# Several match rows had been written in the wrong order (the two rows of a
# same-date pairing were transposed). This swaps the data back (columns B:AC,
# i.e. everything except the running row-index in column A) between each
# affected pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Denmark Division 1")

$rowPairs = @(
    @(280, 281),
    @(337, 338),
    @(362, 363),
    @(371, 372),
    @(439, 440),
    @(540, 541),
    @(547, 548)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Cell-by-cell swap (keeps each value's own type - number/string/date -
    # intact rather than pushing the whole row through a single Variant
    # array, and avoids touching column A which holds the stable row index).
    for ($col = 2; $col -le 29; $col++) {
        $c1 = $ws.Cells.Item($r1, $col)
        $c2 = $ws.Cells.Item($r2, $col)

        $v1 = $c1.Value()
        $v2 = $c2.Value()

        $c1.Value = $v2
        $c2.Value = $v1
    }
}
